$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet 2")

# Add a new sheet after the last existing sheet and name it "Sheet 3"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Sheet 3"

# Copy the header row and the Amazon data row (row 5) from Sheet 2 into the new sheet
$ws2.Range("A1:D1").Copy()
$newSheet.Range("A1").PasteSpecial()
$ws2.Range("A5:D5").Copy()
$newSheet.Range("A2").PasteSpecial()

# Update selections: Sheet 2 selection moves to A5:D5 (no longer the active tab)
$ws2.Range("A5:D5").Select()

# Sheet 3 becomes the active tab, with its data row selected
$newSheet.Range("A2:D2").Select()
